# working on joint mdl test data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: fill in the site/species labels that were missing on rows 135-138
# (they follow the same EM/Noccea pattern as row 134 directly above them)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A135").Value = "EM"
$ws1.Range("B135").Value = "Noccea"
$ws1.Range("A136").Value = "EM"
$ws1.Range("B136").Value = "Noccea"
$ws1.Range("A137").Value = "EM"
$ws1.Range("B137").Value = "Noccea"
$ws1.Range("A138").Value = "EM"
$ws1.Range("B138").Value = "Noccea"

# ---------------------------------------------------------------------------
# Sheet2: fill in missing RM/species labels + G/H (area, dry.mass) columns,
# plus one corrected F value, for the joint mdl test rows (1680-1729)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("F1680").Value = 0.08392

$ws2.Range("A1681").Value = "RM"
$ws2.Range("C1681").Value = "Viola"
$ws2.Range("A1682").Value = "RM"
$ws2.Range("C1682").Value = "Viola"
$ws2.Range("A1683").Value = "RM"
$ws2.Range("C1683").Value = "Viola"
$ws2.Range("A1684").Value = "RM"
$ws2.Range("C1684").Value = "Viola"
$ws2.Range("A1685").Value = "RM"
$ws2.Range("C1685").Value = "Viola"
$ws2.Range("A1686").Value = "RM"
$ws2.Range("C1686").Value = "Viola"
$ws2.Range("A1687").Value = "RM"
$ws2.Range("C1687").Value = "Viola"

$ws2.Range("A1688").Value = "RM"
$ws2.Range("C1688").Value = "Viola"
$ws2.Range("G1688").Value = 4.848
$ws2.Range("H1688").Value = 0.02275

$ws2.Range("A1689").Value = "RM"
$ws2.Range("C1689").Value = "Viola"
$ws2.Range("G1689").Value = 4.271
$ws2.Range("H1689").Value = 0.0197

$ws2.Range("A1690").Value = "RM"
$ws2.Range("C1690").Value = "Androsace"
$ws2.Range("G1690").Value = 0.489
$ws2.Range("H1690").Value = 0.00276

$ws2.Range("A1691").Value = "RM"
$ws2.Range("C1691").Value = "Androsace"
$ws2.Range("G1691").Value = 0.438
$ws2.Range("H1691").Value = 0.00262

$ws2.Range("A1692").Value = "RM"
$ws2.Range("C1692").Value = "Androsace"
$ws2.Range("G1692").Value = 0.916
$ws2.Range("H1692").Value = 0.00412

$ws2.Range("A1693").Value = "RM"
$ws2.Range("C1693").Value = "Androsace"
$ws2.Range("G1693").Value = 0.907
$ws2.Range("H1693").Value = 0.00496

$ws2.Range("A1694").Value = "RM"
$ws2.Range("C1694").Value = "Androsace"
$ws2.Range("G1694").Value = 0.826
$ws2.Range("H1694").Value = 0.00474

$ws2.Range("A1695").Value = "RM"
$ws2.Range("C1695").Value = "Androsace"
$ws2.Range("G1695").Value = 0.859
$ws2.Range("H1695").Value = 0.00483

$ws2.Range("A1696").Value = "RM"
$ws2.Range("C1696").Value = "Androsace"
$ws2.Range("G1696").Value = 0.519
$ws2.Range("H1696").Value = 0.00244

$ws2.Range("A1697").Value = "RM"
$ws2.Range("C1697").Value = "Androsace"
$ws2.Range("G1697").Value = 0.656
$ws2.Range("H1697").Value = 0.00296

$ws2.Range("A1698").Value = "RM"
$ws2.Range("C1698").Value = "Androsace"
$ws2.Range("G1698").Value = 0.71
$ws2.Range("H1698").Value = 0.00271

$ws2.Range("A1699").Value = "RM"
$ws2.Range("C1699").Value = "Androsace"
$ws2.Range("G1699").Value = 0.899
$ws2.Range("H1699").Value = 0.00357

$ws2.Range("A1700").Value = "RM"
$ws2.Range("C1700").Value = "Mertensia f"
$ws2.Range("G1700").Value = 5.19
$ws2.Range("H1700").Value = 0.02561

$ws2.Range("A1701").Value = "RM"
$ws2.Range("C1701").Value = "Mertensia f"
$ws2.Range("G1701").Value = 4.942
$ws2.Range("H1701").Value = 0.02422

$ws2.Range("A1702").Value = "RM"
$ws2.Range("C1702").Value = "Mertensia f"
$ws2.Range("G1702").Value = 2.363
$ws2.Range("H1702").Value = 0.01133

$ws2.Range("A1703").Value = "RM"
$ws2.Range("C1703").Value = "Mertensia f"
$ws2.Range("G1703").Value = 2.439
$ws2.Range("H1703").Value = 0.01094

$ws2.Range("A1704").Value = "RM"
$ws2.Range("C1704").Value = "Mertensia f"
$ws2.Range("G1704").Value = 7.527
$ws2.Range("H1704").Value = 0.0195

$ws2.Range("A1705").Value = "RM"
$ws2.Range("C1705").Value = "Mertensia f"
$ws2.Range("G1705").Value = 3.392
$ws2.Range("H1705").Value = 0.02014

$ws2.Range("A1706").Value = "RM"
$ws2.Range("C1706").Value = "Mertensia f"
$ws2.Range("G1706").Value = 3.236
$ws2.Range("H1706").Value = 0.01296

$ws2.Range("A1707").Value = "RM"
$ws2.Range("C1707").Value = "Mertensia f"
$ws2.Range("G1707").Value = 3.392
$ws2.Range("H1707").Value = 0.01271

$ws2.Range("A1708").Value = "RM"
$ws2.Range("C1708").Value = "Mertensia f"
$ws2.Range("G1708").Value = 5.235
$ws2.Range("H1708").Value = 0.02317

$ws2.Range("A1709").Value = "RM"
$ws2.Range("C1709").Value = "Mertensia f"
$ws2.Range("G1709").Value = 4.642
$ws2.Range("H1709").Value = 0.02179

$ws2.Range("A1710").Value = "RM"
$ws2.Range("F1710").Value = 0.07526
$ws2.Range("G1710").Value = 3.587
$ws2.Range("H1710").Value = 0.0157

$ws2.Range("F1711").Value = 0.08298
$ws2.Range("G1711").Value = 3.917
$ws2.Range("H1711").Value = 0.01603

$ws2.Range("F1712").Value = 0.10551
$ws2.Range("G1712").Value = 5.155
$ws2.Range("H1712").Value = 0.02096

$ws2.Range("F1713").Value = 0.14855
$ws2.Range("G1713").Value = 6.855
$ws2.Range("H1713").Value = 0.02691

$ws2.Range("F1714").Value = 0.03927
$ws2.Range("G1714").Value = 2.109
$ws2.Range("H1714").Value = 0.00796

$ws2.Range("F1715").Value = 0.06709
$ws2.Range("G1715").Value = 3.042
$ws2.Range("H1715").Value = 0.0122

$ws2.Range("F1716").Value = 0.0961
$ws2.Range("G1716").Value = 4.438
$ws2.Range("H1716").Value = 0.01725

$ws2.Range("F1717").Value = 0.17316
$ws2.Range("G1717").Value = 5.995
$ws2.Range("H1717").Value = 0.02931

$ws2.Range("F1718").Value = 0.10885
$ws2.Range("G1718").Value = 4.119
$ws2.Range("H1718").Value = 0.01975

$ws2.Range("F1719").Value = 0.06626
$ws2.Range("G1719").Value = 2.993
$ws2.Range("H1719").Value = 0.01251

$ws2.Range("F1720").Value = 0.06022
$ws2.Range("G1720").Value = 2.766
$ws2.Range("H1720").Value = 0.01193

$ws2.Range("F1721").Value = 0.0737
$ws2.Range("G1721").Value = 3.424
$ws2.Range("H1721").Value = 0.01401

$ws2.Range("F1722").Value = 0.13911
$ws2.Range("G1722").Value = 6.085
$ws2.Range("H1722").Value = 0.02082

$ws2.Range("F1723").Value = 0.11248
$ws2.Range("G1723").Value = 5.198
$ws2.Range("H1723").Value = 0.01894

$ws2.Range("F1724").Value = 0.10992
$ws2.Range("G1724").Value = 4.811
$ws2.Range("H1724").Value = 0.01962

$ws2.Range("F1725").Value = 0.10942
$ws2.Range("G1725").Value = 4.391
$ws2.Range("H1725").Value = 0.01652

$ws2.Range("F1726").Value = 0.07055
$ws2.Range("G1726").Value = 2.818
$ws2.Range("H1726").Value = 0.01323

$ws2.Range("F1727").Value = 0.08172
$ws2.Range("G1727").Value = 3.644
$ws2.Range("H1727").Value = 0.01475

$ws2.Range("F1728").Value = 0.07608
$ws2.Range("G1728").Value = 3.738
$ws2.Range("H1728").Value = 0.01469

$ws2.Range("F1729").Value = 0.05722
$ws2.Range("G1729").Value = 2.936
$ws2.Range("H1729").Value = 0.01102

# ---------------------------------------------------------------------------
# Sheet3: row 24 had header labels mistakenly pasted across columns E, G-L;
# correct E24/H24 and fill in the new G/I/J/K/L label cells.
# Then rename the sheet to "scratch".
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("E24").Value = "site"
$ws3.Range("G24").Value = "species"
$ws3.Range("H24").Value = "indiv"
$ws3.Range("I24").Value = "leaf.no"
$ws3.Range("J24").Value = "fresh.mass"
$ws3.Range("K24").Value = "area"
$ws3.Range("L24").Value = "dry.mass"

$ws3.Name = "scratch"
